$d = $word.ActiveDocument

# Move to the very end of the document (after the last paragraph,
# "O JUnit Teste resulta em not null usando o enum Resposta no método
# calculaPena.") and add two new paragraphs:
#   1) an empty paragraph
#   2) a paragraph with the new report text
# Typing it as two separate paragraph breaks (with nothing typed in
# between the first break and the second) reproduces exactly that:
# the first new paragraph stays empty, the second receives the text.

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Content
$r.Collapse(0)
$r.InsertAfter("O JUnit Teste na terceira rodada resulta em erro, não foi possível usar os valores nas variáveis da classe JulgamentoPrisioneiro, os getters e setters devem ser criados")
